$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly price observation was added for "Papa" at Terminal
# Hortofrutícola Agro Chillán. It belongs at the top of this block (row 179),
# so insert a fresh row there and push the existing rows 179-194 down to 180-195.
$ws.Rows.Item(179).Insert()

# Fill in the new row 179 with the new observation's data.
$ws.Range("A179").Value = 7
$ws.Range("B179").Value = 'Terminal Hortofrutícola Agro Chillán'
$ws.Range("C179").Value = 'Ñuble'
$ws.Range("D179").Value = 44461
$ws.Range("E179").Value = 16
$ws.Range("F179").Value = 100114001
$ws.Range("G179").Value = 'Papa'
$ws.Range("H179").Value = 'Patagonia'
$ws.Range("I179").Value = '1a (guarda)'
$ws.Range("J179").Value = 300
$ws.Range("K179").Value = 7000
$ws.Range("L179").Value = 7500
$ws.Range("M179").Value = 7250
$ws.Range("N179").Value = '$/saco 25 kilos'
$ws.Range("O179").Value = 'Provincia de Diguillín'
$ws.Range("P179").Value = 290
$ws.Range("Q179").Value = 25
$ws.Range("R179").Value = 'Hortaliza'
